$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data: replace the PACS registration record with the new one ---
$ws.Range("D2").Value = "Jharsuguda"
$ws.Range("E2").Value = "Jharsuguda"
$ws.Range("B2").Value = "sasi11@gmail.com"
$ws.Range("C2").Value = "Siba@123"
$ws.Range("F2").Value = "Laikera"
$ws.Range("G2").Value = "LAIKERA SCS"
$ws.Range("A2").Value = "Sasi"

# Name Of PPC (I2) stays "Rama" - no content change needed.

# PACS Submit status (O2) was cleared out entirely.
$ws.Range("O2").ClearContents()

# New hyperlink on the PACS Password cell (C2), mirroring the existing
# PACS EmailId (B2) / ARCS EmailId (Q2) mailto-style hyperlinks.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Siba@123")

# --- Column D widened to fit the new "Jharsuguda" / header text ---
$ws.Columns.Item(4).ColumnWidth = 25

# --- Selection moved to A4, with no frozen/scrolled top-left cell ---
$ws.Range("A4").Select()
